$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (columns C/D/E relabeled; A/B unchanged) ---
$ws.Range("C1").Value = "sd log10 cas3/gapA ratio"
$ws.Range("D1").Value = "mean log10 cas3/gapA ratio"
$ws.Range("E1").Value = "n"

# --- Replace the data block (rows 2-13) with the new measurements ---
$data = @(
    @("α15.2",  5, 0.0855699872200405, 0.190835096524283, 6),
    @("α15.2", 15, 0.372258295505798,  0.878060865572723, 6),
    @("α15.2", 30, 0.381402668197424,  1.5478036469452,   6),
    @("α20.4",  5, 0.0773842607499961, 0.191166591886716, 6),
    @("α20.4", 15, 1.38032687535014,   1.80923257560708,  6),
    @("α20.4", 30, 0.389084056008771,  1.49260992397535,  6),
    @("α48.4",  5, 0.0166862277138728, 0.178920610897078, 6),
    @("α48.4", 15, 0.248178000509146,  0.626128902307533, 6),
    @("α48.4", 30, 0.490297055221236,  3.10016104803487,  6),
    @("α51.5",  5, 0.0149318666345928, 0.0909223108225965,6),
    @("α51.5", 15, 0.0843655926173616, 0.298269280049384, 6),
    @("α51.5", 30, 0.545897142407056,  2.18403319359579,  6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}

# --- The old sheet had 16 data rows (2-17); only 12 remain now, so drop the rest ---
$ws.Rows("14:17").Delete()
